# Applies the S23/G06 "Alerts refactor v3" sprint-task updates:
#  - bump row heights for rows 195-198 (minor re-wrap height tweaks)
#  - append rows 199-212 documenting the new S23_G06 backend/frontend tasks
#  - update the sheet view (active cell) to reflect the new bottom of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row-height tweaks on existing rows 195-198 (content unchanged) ---
$ws.Rows.Item(195).RowHeight = 41.25
$ws.Rows.Item(196).RowHeight = 41.25
$ws.Rows.Item(197).RowHeight = 27.75
$ws.Rows.Item(198).RowHeight = 41.25

# --- New rows 199-212: S23 / G06 "Alerts refactor v3" tasks ---
# Row 199 (S23_G06_TB001)
$ws.Range("A199").Value = 'S23'
$ws.Range("B199").Value = 'G06'
$ws.Range("C199").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D199").Value = 'S23_G06_TB001'
$ws.Range("E199").Value = 'Define AlertDefinition + AlertEvent schemas/models (target symbol/universe, variables, condition DSL, trigger_mode, evaluation_cadence, time constraints) and CRUD APIs.'
$ws.Range("F199").Value = 'New v3 tables + APIs under /api/alerts-v3; alert targets: SYMBOL/HOLDINGS/GROUP.'
$ws.Range("F199").Style = "Normal"
$ws.Range("G199").Value = 'implemented'
$ws.Range("H199").Value = 'Added AlertDefinition/CustomIndicator/AlertEvent models + CRUD.'
$ws.Rows.Item(199).RowHeight = 55.2

# Row 200 (S23_G06_TB002)
$ws.Range("A200").Value = 'S23'
$ws.Range("B200").Value = 'G06'
$ws.Range("C200").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D200").Value = 'S23_G06_TB002'
$ws.Range("E200").Value = 'Implement per-alert evaluation cadence scheduler using latest completed bars for each referenced timeframe; enforce missing-data=false rule.'
$ws.Range("F200").Value = 'Scheduler runs every ~15s and skips alerts until cadence due; uses latest available bars in DB.'
$ws.Range("F200").Style = "Normal"
$ws.Range("G200").Value = 'implemented'
$ws.Range("H200").Value = 'Per-alert cadence evaluation implemented; missing-data returns false.'
$ws.Rows.Item(200).RowHeight = 41.75

# Row 201 (S23_G06_TB003)
$ws.Range("A201").Value = 'S23'
$ws.Range("B201").Value = 'G06'
$ws.Range("C201").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D201").Value = 'S23_G06_TB003'
$ws.Range("E201").Value = 'Implement event operators semantics: CROSSES_ABOVE/BELOW and MOVING_UP/DOWN (numeric RHS only) with now/prev rules.'
$ws.Range("F201").Value = 'Supports aliases CROSSING_* -> CROSSES_*.'
$ws.Range("F201").Style = "Normal"
$ws.Range("G201").Value = 'implemented'
$ws.Range("H201").Value = 'CROSSES_ABOVE/BELOW + MOVING_UP/DOWN (numeric RHS) implemented.'
$ws.Rows.Item(201).RowHeight = 41.75

# Row 202 (S23_G06_TB004)
$ws.Range("A202").Value = 'S23'
$ws.Range("B202").Value = 'G06'
$ws.Range("C202").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D202").Value = 'S23_G06_TB004'
$ws.Range("E202").Value = 'Add metrics/columns as operands (TODAY_PNL_PCT, PNL_PCT, INVESTED, CURRENT_VALUE, etc.) accessible in expressions and snapshots.'
$ws.Range("F202").Value = 'Metrics computed from Position + 1d candles; prev values approximated.'
$ws.Range("F202").Style = "Normal"
$ws.Range("G202").Value = 'implemented'
$ws.Range("H202").Value = 'Added metric operands (TODAY_PNL_PCT, PNL_PCT, etc).'
$ws.Range("I202").Value = 'Enhance snapshot to include per-variable values.'
$ws.Range("I202").Style = "Normal"
$ws.Rows.Item(202).RowHeight = 41.75

# Row 203 (S23_G06_TB005)
$ws.Range("A203").Value = 'S23'
$ws.Range("B203").Value = 'G06'
$ws.Range("C203").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D203").Value = 'S23_G06_TB005'
$ws.Range("E203").Value = 'Custom indicators backend (Phase A): model + CRUD + validation + allowed function set (A MVP surface) + compilation/cache hooks.'
$ws.Range("F203").Value = 'Phase A function allowlist enforced at compile time; recursion not supported.'
$ws.Range("F203").Style = "Normal"
$ws.Range("G203").Value = 'implemented'
$ws.Range("H203").Value = 'Custom indicators CRUD + validation + compilation implemented.'
$ws.Range("I203").Value = 'Add preview endpoint/UI for formula values.'
$ws.Range("I203").Style = "Normal"
$ws.Rows.Item(203).RowHeight = 41.75

# Row 204 (S23_G06_TB006)
$ws.Range("A204").Value = 'S23'
$ws.Range("B204").Value = 'G06'
$ws.Range("C204").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D204").Value = 'S23_G06_TB006'
$ws.Range("E204").Value = 'Add evaluation/test endpoints for “Test on last bar” preview (return per-symbol TRUE/FALSE + snapshot + missing-data reasons).'
$ws.Range("G204").Value = 'pending'
$ws.Range("H204").Value = 'Not implemented yet.'
$ws.Range("I204").Value = 'Add test/preview endpoint returning per-symbol evaluation + snapshot/missing-data reason.'
$ws.Range("I204").Style = "Normal"
$ws.Rows.Item(204).RowHeight = 41.75

# Row 205 (S23_G06_TB007)
$ws.Range("A205").Value = 'S23'
$ws.Range("B205").Value = 'G06'
$ws.Range("C205").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D205").Value = 'S23_G06_TB007'
$ws.Range("E205").Value = 'Backend tests: condition builder DSL serialization, operator semantics, per-alert cadence scheduling, custom indicator validation guardrails.'
$ws.Range("G205").Value = 'implemented'
$ws.Range("H205").Value = 'Added backend regression tests for v3 parser/compiler/evaluator + API.'
$ws.Rows.Item(205).RowHeight = 41.75

# Row 206 (S23_G06_TF001)
$ws.Range("A206").Value = 'S23'
$ws.Range("B206").Value = 'G06'
$ws.Range("C206").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D206").Value = 'S23_G06_TF001'
$ws.Range("E206").Value = 'Build Alerts module page with tabs: Alerts / Indicators / Events; add routing + navigation entry.'
$ws.Range("F206").Value = 'Implemented inside existing Alerts page as tabs; kept Legacy tab.'
$ws.Range("F206").Style = "Normal"
$ws.Range("G206").Value = 'implemented'
$ws.Range("H206").Value = 'Alerts page now has tabs: Alerts/Indicators/Events/Legacy.'
$ws.Rows.Item(206).RowHeight = 28.35

# Row 207 (S23_G06_TF002)
$ws.Range("A207").Value = 'S23'
$ws.Range("B207").Value = 'G06'
$ws.Range("C207").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D207").Value = 'S23_G06_TF002'
$ws.Range("E207").Value = 'Implement Create/Edit Alert wizard: Target → Variables → Conditions → Trigger settings; store canonical DSL string; show read-only DSL preview.'
$ws.Range("F207").Value = 'Implemented as a single Create/Edit dialog (not multi-step wizard).'
$ws.Range("F207").Style = "Normal"
$ws.Range("G207").Value = 'implemented'
$ws.Range("H207").Value = 'Create/Edit alert via dialog; variables defined as name+DSL.'
$ws.Range("I207").Value = 'Add multi-step wizard + target/group picker UI.'
$ws.Range("I207").Style = "Normal"
$ws.Rows.Item(207).RowHeight = 41.75

# Row 208 (S23_G06_TF003)
$ws.Range("A208").Value = 'S23'
$ws.Range("B208").Value = 'G06'
$ws.Range("C208").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D208").Value = 'S23_G06_TF003'
$ws.Range("E208").Value = 'Condition builder UI: operand pickers (Variable/Metric/Constant), operators (relational + event), AND/OR join, add/remove rows.'
$ws.Range("G208").Value = 'planned'
$ws.Range("H208").Value = 'Not implemented; using free-form DSL for conditions.'
$ws.Range("I208").Value = 'Build condition builder UI with operand pickers and AND/OR joining.'
$ws.Range("I208").Style = "Normal"
$ws.Rows.Item(208).RowHeight = 41.75

# Row 209 (S23_G06_TF004)
$ws.Range("A209").Value = 'S23'
$ws.Range("B209").Value = 'G06'
$ws.Range("C209").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D209").Value = 'S23_G06_TF004'
$ws.Range("E209").Value = 'Variables UI: indicator variable rows + metric variable support; hide/disable Bars/Length for PRICE/VOLUME primitives.'
$ws.Range("F209").Value = 'Variable UI currently DSL-only; no structured indicator/metric variable rows yet.'
$ws.Range("F209").Style = "Normal"
$ws.Range("G209").Value = 'planned'
$ws.Range("H209").Value = 'Not implemented.'
$ws.Range("I209").Value = 'Add structured variable builder (indicator/metric) + hide non-applicable fields.'
$ws.Range("I209").Style = "Normal"
$ws.Rows.Item(209).RowHeight = 41.75

# Row 210 (S23_G06_TF005)
$ws.Range("A210").Value = 'S23'
$ws.Range("B210").Value = 'G06'
$ws.Range("C210").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D210").Value = 'S23_G06_TF005'
$ws.Range("E210").Value = 'Indicators tab UI: custom indicator list + create/edit dialog with formula builder and allowed-function guidance (Phase A).'
$ws.Range("F210").Value = 'Formula is free-form DSL textarea; no interactive builder/preview.'
$ws.Range("F210").Style = "Normal"
$ws.Range("G210").Value = 'implemented'
$ws.Range("H210").Value = 'Indicators tab supports custom indicator CRUD.'
$ws.Range("I210").Value = 'Add guided builder + preview.'
$ws.Range("I210").Style = "Normal"
$ws.Rows.Item(210).RowHeight = 41.75

# Row 211 (S23_G06_TF006)
$ws.Range("A211").Value = 'S23'
$ws.Range("B211").Value = 'G06'
$ws.Range("C211").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D211").Value = 'S23_G06_TF006'
$ws.Range("E211").Value = 'Events tab UI: list/filter alert events and display snapshots/reasons; link back to alert definition.'
$ws.Range("F211").Value = 'Events tab lists events; snapshot drilldown not shown yet.'
$ws.Range("F211").Style = "Normal"
$ws.Range("G211").Value = 'implemented'
$ws.Range("H211").Value = 'Events tab implemented (basic list).'
$ws.Range("I211").Value = 'Show snapshot JSON + link back to alert.'
$ws.Range("I211").Style = "Normal"
$ws.Rows.Item(211).RowHeight = 28.35

# Row 212 (S23_G06_TF007)
$ws.Range("A212").Value = 'S23'
$ws.Range("B212").Value = 'G06'
$ws.Range("C212").Value = 'Alerts refactor v3: indicator-first alerts over universes (see docs/alerts_refactor_v3.md)'
$ws.Range("D212").Value = 'S23_G06_TF007'
$ws.Range("E212").Value = 'Frontend tests: wizard flow smoke test + condition builder serialization + selected-rows→group guidance.'
$ws.Range("G212").Value = 'pending'
$ws.Range("H212").Value = 'Not implemented.'
$ws.Range("I212").Value = 'Add minimal frontend tests for v3 alerts UI.'
$ws.Range("I212").Style = "Normal"
$ws.Rows.Item(212).RowHeight = 28.35

# --- Sheet view: scroll target / active selection moved to the new bottom rows ---
[void]$ws.Range("C207").Select()

Write-Output "Applied S23/G06 Alerts refactor v3 rows (199-212) and row-height updates."
